# Fruta / hortaliza, semanal
# Insert two new weekly rows at the top of the data block (rows 6-7),
# pushing the existing data down by 2 rows (old row 6 -> row 8, ...,
# old row 73 -> row 75), and populate the two new rows with this week's
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 6.. down by two rows so the new week's rows land at 6 and 7.
$ws.Range("A6:R7").Insert()

# New row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 45190
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 100112028
$ws.Range("G6").Value = "Sandia"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 650
$ws.Range("K6").Value = 700
$ws.Range("L6").Value = 730
$ws.Range("M6").Value = 716
$ws.Range("N6").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 716
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = "Hortaliza"

# New row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 45190
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112028
$ws.Range("G7").Value = "Sandia"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Tercera"
$ws.Range("J7").Value = 710
$ws.Range("K7").Value = 700
$ws.Range("L7").Value = 730
$ws.Range("M7").Value = 714
$ws.Range("N7").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 714
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
